$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "20/12/2023 16:05:59"
$ws.Range("E3").Value = "20/12/2023 16:05:59"
$ws.Range("E4").Value = "20/12/2023 16:05:59"
$ws.Range("E5").Value = "20/12/2023 16:05:59"
$ws.Range("E6").Value = "20/12/2023 16:05:59"
$ws.Range("E7").Value = "20/12/2023 16:05:59"
$ws.Range("E8").Value = "20/12/2023 16:05:59"
$ws.Range("E9").Value = "20/12/2023 16:05:59"
$ws.Range("E10").Value = "20/12/2023 16:05:59"
$ws.Range("E11").Value = "20/12/2023 16:05:59"
$ws.Range("E12").Value = "20/12/2023 16:05:59"
$ws.Range("E13").Value = "20/12/2023 16:05:59"

$ws.Range("C4").Value = "Ex de mulher sequestrada com Marcelinho teme ser atacado na rua"

$ws.Range("C11").Value = "Lutador do UFC agradece rival por parar golpes: “Poderia ter me matado”"
$ws.Range("D11").Value = "https://www.cnnbrasil.com.br/noticias/lutador-do-ufc-agradece-rival-por-parar-golpes-poderia-ter-me-matado/"

$ws.Range("C12").Value = "Estados desistem de aumentar ICMS após mudança na reforma tributária"
$ws.Range("D12").Value = "https://www.cnnbrasil.com.br/economia/estados-desistem-de-aumentar-icms-apos-mudanca-na-reforma-tributaria/"

$ws.Range("C13").Value = "Vendas do PS5 ultrapassam 50 milhões de unidades, diz Sony"
$ws.Range("D13").Value = "https://www.cnnbrasil.com.br/economia/vendas-do-ps5-ultrapassam-50-milhoes-de-unidades-diz-sony/"
